# Apply the cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "46.972.11"
$ws.Range("E2").Value = "  +5.73%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.336.64"
$ws.Range("E3").Value = "  +5.14%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.75%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "305.59"
$ws.Range("E5").Value = "  +0.94%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "97.50"
$ws.Range("E6").Value = "  +8.31%  "

# --- Row 7: XRP ---
$ws.Range("E7").Value = "  +3.98%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  -0.69%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +8.37%  "

# --- Row 10: Avalanche ---
$ws.Range("D10").Value = "35.99"
$ws.Range("E10").Value = "  +6.73%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  +3.92%  "

# --- Row 12: Polkadot ---
$ws.Range("E12").Value = "  +7.28%  "

# --- Row 13: TRON ---
$ws.Range("E13").Value = "  +0.09%  "

# --- Row 14: WrappedliquidstakedEther2.0 ---
$ws.Range("D14").Value = "2.691.98"
$ws.Range("E14").Value = "  +4.89%  "

# --- Row 15: WrappedEther ---
$ws.Range("D15").Value = "2.339.30"
$ws.Range("E15").Value = "  +0.35%  "

# --- Row 16: Chainlink ---
$ws.Range("E16").Value = "  +7.61%  "

# --- Row 17: Polygon ---
$ws.Range("D17").Value = "0.838"
$ws.Range("E17").Value = "  +3.87%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "46.820.92"
$ws.Range("E18").Value = "  +5.92%  "

# --- Row 19: InternetComputer(DFINITY) ---
$ws.Range("D19").Value = "13.74"
$ws.Range("E19").Value = "  +20.73%  "

# --- Row 20: ShibaInu ---
$ws.Range("E20").Value = "  +4.45%  "

# --- Row 21: Uniswap ---
$ws.Range("E21").Value = "  +3.17%  "

# --- Row 22: Litecoin ---
$ws.Range("D22").Value = "68.03"
$ws.Range("E22").Value = "  +5.66%  "

# --- Row 23: BitcoinCash ---
$ws.Range("D23").Value = "251.68"
$ws.Range("E23").Value = "  +7.65%  "

# --- Row 24: PancakeSwap ---
$ws.Range("E24").Value = "  +3.95%  "

# --- Row 25: ImmutableX ---
$ws.Range("E25").Value = "  +4.38%  "

# --- Row 26: Dai ---
$ws.Range("E26").Value = "  -0.32%  "

# --- Row 27: InjectiveProtocol ---
$ws.Range("D27").Value = "42.57"
$ws.Range("E27").Value = "  +17.16%  "

# --- Row 28: Toncoin ---
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").Value = "  +1.64%  "

# --- Row 29: Cosmos ---
$ws.Range("D29").Value = "9.93"
$ws.Range("E29").Value = "  +4.59%  "

# --- Row 30: EthereumClassic ---
$ws.Range("D30").Value = "20.35"
$ws.Range("E30").Value = "  +4.07%  "

# --- Row 31: Filecoin ---
$ws.Range("D31").Value = "5.86"
$ws.Range("E31").Value = "  +4.42%  "

# --- Row 32: Hedera ---
$ws.Range("E32").Value = "  +7.77%  "

# --- Row 33: Monero ---
$ws.Range("D33").Value = "146.86"
$ws.Range("E33").Value = "  +0.55%  "

# --- Row 34: WEMIXToken ---
$ws.Range("E34").Value = "  -0.09%  "

# --- Rows 35 & 36: Kaspa and LidoDAOToken swap positions and update values ---
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "3.15"
$ws.Range("E35").Value = "  +5.00%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.115"
$ws.Range("E36").Value = "  +7.86%  "

# --- Row 37: Stellar ---
$ws.Range("E37").Value = "  +3.58%  "

# --- Row 38: ARBITRUM ---
$ws.Range("E38").Value = "  +2.21%  "

# --- Row 39: RenderToken ---
$ws.Range("D39").Value = "4.01"
$ws.Range("E39").Value = "  +9.66%  "

# --- Row 40: VeChain ---
$ws.Range("E40").Value = "  +8.05%  "

# --- Row 41: NEARProtocol ---
$ws.Range("E41").Value = "  +5.38%  "

# --- Row 42: Celestia ---
$ws.Range("D42").Value = "14.04"
$ws.Range("E42").Value = "  -4.37%  "

# --- Row 43: FirstDigitalUSD ---
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.94%  "

# --- Row 44: Stacks ---
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +14.90%  "

# --- Row 45: Maker ---
$ws.Range("D45").Value = "1.804.58"
$ws.Range("E45").Value = "  +1.75%  "

# --- Row 46: BitcoinSV ---
$ws.Range("D46").Value = "91.63"
$ws.Range("E46").Value = "  +15.68%  "

# --- Row 47: ordi ---
$ws.Range("D47").Value = "75.03"
$ws.Range("E47").Value = "  +11.88%  "

# --- Row 48: Algorand ---
$ws.Range("E48").Value = "  +7.71%  "

# --- Row 49: Aave ---
$ws.Range("D49").Value = "99.01"
$ws.Range("E49").Value = "  +3.79%  "

# --- Row 50: MultiversX ---
$ws.Range("D50").Value = "55.43"
$ws.Range("E50").Value = "  +5.84%  "

# --- Row 51: FraxShare ---
$ws.Range("E51").Value = "  +5.75%  "
